$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data inserted at the top (row 2), with all subsequent rows
# (previously rows 2-10) shifted down by one row. The last previous row
# (row 11) is pushed out of the table range and dropped. Column A (the
# Q-index labels) stays fixed per row.

$newValues = @(
    @(-0.02314597604078636, 0.3579920056255013, 0.1782699060034266, 0.4222202103209018, 0.4363822494547141, 15),
    @(0.01177902937386835, 0.3233677337305382, 0.1387525778210147, 0.3724950708680783, 0.3863630867095333, 14),
    @(-0.008234150994489252, 0.3656812685479343, 0.1736781030421116, 0.4167470492302394, 0.4336794057244024, 13),
    @(0.04114622104959433, 0.2438785208639479, 0.09301430872978532, 0.3049824728239072, 0.315631483433317, 12),
    @(-0.004251297071472607, 0.3196557001756482, 0.1664324393548997, 0.4079613209054257, 0.4278502102499434, 11),
    @(-0.03152289704198964, 0.3137577720403564, 0.1200567451426497, 0.3464920563918454, 0.3637200573468097, 10),
    @(-0.05669611638992965, 0.309241221828214, 0.12362493933599, 0.3516033835673229, 0.3680513563860018, 9),
    @(-0.04410238231315148, 0.2487243691024244, 0.1063280554679864, 0.3260798299005727, 0.3453909343541705, 8),
    @(-0.07289500678705241, 0.3460499651212735, 0.2127274776992719, 0.4612238910759848, 0.4919174337551497, 7),
    @(-0.05253816173664939, 0.2093953811075043, 0.05818608614024737, 0.2412179225104291, 0.2578972532655152, 6)
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $vals = $newValues[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
